$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (J1:P1) --------------------------------------------
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "renterpov"
$ws.Range("L1").Value = "ownerpov"
$ws.Range("M1").Value = "burden50r"
$ws.Range("N1").Value = "burden50h"
$ws.Range("O1").Value = "burden30r"
$ws.Range("P1").Value = "burden30h"

# --- 2005 row -------------------------------------------------------------
$ws.Range("J2").Value = 2005
$ws.Range("K2").Formula = "=D2"
$ws.Range("L2").Formula = "=D13"
$ws.Range("M2").Formula = "=E7"
$ws.Range("N2").Formula = "=E18"
$ws.Range("O2").Formula = "=D7"
$ws.Range("P2").Formula = "=D18"

# --- 2010 row -------------------------------------------------------------
$ws.Range("J3").Value = 2010
$ws.Range("K3").Formula = "=D3"
$ws.Range("L3").Formula = "=D14"
$ws.Range("M3").Formula = "=E8"
$ws.Range("N3").Formula = "=E19"
$ws.Range("O3").Formula = "=D8"
$ws.Range("P3").Formula = "=D19"

# --- 2015 row -------------------------------------------------------------
$ws.Range("J4").Value = 2015
$ws.Range("K4").Formula = "=D4"
$ws.Range("L4").Formula = "=D15"
$ws.Range("M4").Formula = "=E9"
$ws.Range("N4").Formula = "=E20"
$ws.Range("O4").Formula = "=D9"
$ws.Range("P4").Formula = "=D20"

# --- selection moves to the newly-prepared chart-source range -------------
$ws.Range("L1:P4").Select()
